$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.956.74"
$ws.Range("E2").Value = "  +2.35%  "

$ws.Range("D3").Value = "3.009.56"

$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

$c = $ws.Range("D5")
$c.Value = "'515.07"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.90%  "

$c = $ws.Range("D6")
$c.Value = "'139.58"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.15%  "

$c = $ws.Range("D7")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$c = $ws.Range("D8")
$c.Value = "'0.436"
$c.Style = "Normal"

$ws.Range("E9").Value = "  +5.85%  "

$c = $ws.Range("D10")
$c.Value = "'0.109"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +7.31%  "

$c = $ws.Range("D11")
$c.Value = "'0.358"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.26%  "

$ws.Range("E12").Value = "  +2.64%  "

$ws.Range("D13").Value = "3.518.50"
$ws.Range("E13").Value = "  +1.89%  "

$c = $ws.Range("D14")
$c.Value = "'25.77"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.55%  "

$ws.Range("E15").Value = "  +11.95%  "

$ws.Range("D16").Value = "56.946.62"
$ws.Range("E16").Value = "  +2.48%  "

$ws.Range("D17").Value = "3.003.66"
$ws.Range("E17").Value = "  +1.75%  "

$c = $ws.Range("D18")
$c.Value = "'5.96"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.58%  "

$c = $ws.Range("D19")
$c.Value = "'12.61"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.91%  "

$c = $ws.Range("D20")
$c.Value = "'7.89"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.29%  "

$c = $ws.Range("D21")
$c.Value = "'328.42"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.48%  "

$ws.Range("E22").Value = "  -0.09%  "

$c = $ws.Range("D23")
$c.Value = "'0.486"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +5.34%  "

$c = $ws.Range("D24")
$c.Value = "'63.49"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +5.70%  "

$ws.Range("E25").Value = "  +7.08%  "

$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").Value = "0.0₃0921"
$ws.Range("E27").Value = "  +8.91%  "

$c = $ws.Range("D28")
$c.Value = "'6.68"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.43%  "

$c = $ws.Range("D29")
$c.Value = "'7.11"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +8.51%  "

$c = $ws.Range("D30")
$c.Value = "'1.24"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +7.17%  "

$ws.Range("E31").Value = "  +6.91%  "

$c = $ws.Range("D32")
$c.Value = "'20.67"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +6.15%  "

$c = $ws.Range("D33")
$c.Value = "'157.09"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.32%  "

$c = $ws.Range("D34")
$c.Value = "'4.60"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.04%  "

$c = $ws.Range("D35")
$c.Value = "'5.73"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D37")
$c.Value = "'24.39"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.55%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D38")
$c.Value = "'0.0681"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.79%  "

$ws.Range("D39").Value = "3.036.98"
$ws.Range("E39").Value = "  +1.92%  "

$c = $ws.Range("D40")
$c.Value = "'37.19"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "

$c = $ws.Range("D41")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "2.291.42"
$ws.Range("E42").Value = "  +8.43%  "

$c = $ws.Range("D43")
$c.Value = "'0.650"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.51%  "

$c = $ws.Range("D44")
$c.Value = "'1.43"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.24%  "

$c = $ws.Range("D45")
$c.Value = "'3.70"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.64%  "

$c = $ws.Range("D46")
$c.Value = "'1.01"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "

$c = $ws.Range("D47")
$c.Value = "'1.98"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +9.32%  "

$c = $ws.Range("D48")
$c.Value = "'0.0241"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.00%  "

$c = $ws.Range("D49")
$c.Value = "'5.89"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.39%  "

$c = $ws.Range("D50")
$c.Value = "'19.30"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "

$c = $ws.Range("D51")
$c.Value = "'0.0878"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.83%  "
